$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LJ Speech")

# Row 2
$ws.Range("B2").Value = "<kere>"
$ws.Range("C2").Value = 21

# Row 3
$ws.Range("B3").Value = "<antire>"
$ws.Range("C3").Value = 29

# Row 4
$ws.Range("C4").Value = 33

# Row 5
$ws.Range("C5").Value = 39

# Row 6
$ws.Range("C6").Value = 34

# Row 7
$ws.Range("C7").Value = 37

# Row 8
$ws.Range("C8").Value = 42

# Row 9
$ws.Range("B9").Value = "<then>"

# Row 10
$ws.Range("C10").Value = 36

# Row 11
$ws.Range("B11").Value = "<word>"
$ws.Range("C11").Value = 31

# Row 12
$ws.Range("C12").Value = 36

# Row 13
$ws.Range("C13").Value = 32

# Row 14
$ws.Range("C14").Value = 35

# Row 15
$ws.Range("B15").Value = "<cartine>"
$ws.Range("C15").Value = 13
